# Idol.xlsx: rename Min/Max to CurMin/CurMax, update their values, and add
# two new columns PotenMin / PotenMax (E, F) with header + type rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("C1").Value = "CurMin"
$ws.Range("D1").Value = "CurMax"
$ws.Range("E1").Value = "PotenMin"
$ws.Range("F1").Value = "PotenMax"

# --- Type row (row 2) ---
$ws.Range("E2").Value = "int"
$ws.Range("F2").Value = "int"

# --- Data rows (row 3..6) ---
# Row 3 (ID 1, normal)
$ws.Range("D3").Value = 40
$ws.Range("E3").Value = 20
$ws.Range("F3").Value = 70

# Row 4 (ID 2, silver)
$ws.Range("D4").Value = 70
$ws.Range("E4").Value = 40
$ws.Range("F4").Value = 80

# Row 5 (ID 3, gold)
$ws.Range("D5").Value = 80
$ws.Range("E5").Value = 60
$ws.Range("F5").Value = 90

# Row 6 (ID 4, platinum)
$ws.Range("D6").Value = 95
$ws.Range("E6").Value = 80
$ws.Range("F6").Value = 99

# --- Selection matching the post-edit state ---
$ws.Range("G2").Select()

# --- Page setup (adds pageSetup element with paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
